$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refine the "type" tokens for several columns ---------------
# The shared-string table must grow in a very specific order so the
# resulting indices line up with the target workbook. Writing the cells
# below in this exact sequence makes each *new* unique string get
# appended to xl/sharedStrings.xml in the right spot.

# I2 keeps the #float token but now documents the unit (ml or mg).
$ws.Range("I2").Value = "#float,  unit:mlormg"

# J2 (Wavelength) becomes an integer with a nm unit.
$ws.Range("J2").Value = "#integer,  unit:nm"

# M2/N2/P2 (ReagentBlankValue, SampleBlankValue, Result) are plain floats.
$ws.Range("M2").Value = "#float"
$ws.Range("N2").Value = "#float"
$ws.Range("P2").Value = "#float"

# O2 (MolarExtinctionCoefficient) is a float with l/mol/cm unit.
$ws.Range("O2").Value = "#float,  unit:l/mol/cm"

# Q2 (Unit) becomes a string enumerating the allowed unit tokens.
$ws.Range("Q2").Value = "#string,  unit:mmol/lormg/lorµmol/lorg/l"

# --- Row 3: new description row for the headers -------------------------
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"

# The remaining description cells (H3:R3) are blank in the target sheet.
$ws.Range("H3:R3").Value = ""
